$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "B2" = 0.864852945860008
    "C2" = 0.08470519802246201
    "D2" = 0.07666607668650727
    "E2" = 0.1006978799920759
    "G2" = 3.102342321812586
    "H2" = 2.335905672660203
    "K2" = 0.4411727187139434
    "L2" = 0.1966204710330004
    "M2" = 0.2181622898218905
    "B3" = 0.843234501217438
    "C3" = 0.08253124014887803
    "D3" = 0.0698194154473839
    "E3" = 0.1003567323239896
    "G3" = 3.027054271993848
    "H3" = 2.303762807196648
    "K3" = 0.4202695279218602
    "L3" = 0.1943271504768944
    "M3" = 0.2139217819283559
    "B4" = 0.8306423692436908
    "C4" = 0.08115432108868248
    "D4" = 0.06565449457212935
    "E4" = 0.1001677995846375
    "G4" = 2.981526307777557
    "H4" = 2.284513664771538
    "K4" = 0.407858256850389
    "L4" = 0.1930160115647368
    "M4" = 0.2114682761707911
    "B5" = 0.8256823377328431
    "C5" = 0.080582543639796
    "D5" = 0.06396696806146451
    "E5" = 0.1000959758612865
    "G5" = 2.963148183260955
    "H5" = 2.276791536964026
    "K5" = 0.4029068110391876
    "L5" = 0.192506104955612
    "M5" = 0.2105062181628874
    "B6" = 0.8248690770264773
    "C6" = 0.08048695338806056
    "D6" = 0.06368733976729857
    "E6" = 0.1000843618546128
    "G6" = 2.960107045642445
    "H6" = 2.275516645662179
    "K6" = 0.4020910382864571
    "L6" = 0.1924229089363294
    "M6" = 0.2103487501695724
    "B7" = 0.8305747827378411
    "C7" = 0.08114665321979331
    "D7" = 0.06563169678200609
    "E7" = 0.1001668100135991
    "G7" = 2.981277747074245
    "H7" = 2.284409027721267
    "K7" = 0.4077910500042208
    "L7" = 0.1930090360048169
    "M7" = 0.2114551486069942
    "B8" = 0.8572573811734401
    "C8" = 0.0839643035435067
    "D8" = 0.07429720983748211
    "E8" = 0.1005759931287571
    "G8" = 3.076237296896977
    "H8" = 2.324721583962145
    "K8" = 0.4338772895606837
    "L8" = 0.195809612392317
    "M8" = 0.2166689900694898
    "B9" = 0.9149977110835721
    "C9" = 0.0891596199560496
    "D9" = 0.09160405130425886
    "E9" = 0.1015412160338354
    "G9" = 3.268053185758902
    "H9" = 2.407656317084246
    "K9" = 0.4884050215825368
    "L9" = 0.2020710905084684
    "M9" = 0.2280859067562311
    "B10" = 0.9607377671512154
    "C10" = 0.09278102463692761
    "D10" = 0.1045185199519238
    "E10" = 0.1023496232834411
    "G10" = 3.412482847937099
    "H10" = 2.470990281754325
    "K10" = 0.5305467570593407
    "L10" = 0.2071415803825261
    "M10" = 0.2372036934400015
    "B11" = 0.9822707641607167
    "C11" = 0.09438724555126754
    "D11" = 0.1104386625195701
    "E11" = 0.102738953130455
    "G11" = 3.47896818788638
    "H11" = 2.500332168612545
    "K11" = 0.5501754617798724
    "L11" = 0.2095506630355715
    "M11" = 0.2415107815971069
    "B12" = 0.9905293070940786
    "C12" = 0.09498965339215459
    "D12" = 0.1126870942686935
    "E12" = 0.1028894838412455
    "G12" = 3.50425842031882
    "H12" = 2.511520029705821
    "K12" = 0.5576745895192801
    "L12" = 0.2104776668158053
    "M12" = 0.2431647120785669
    "B13" = 0.9887460340019061
    "C13" = 0.09486017211776243
    "D13" = 0.1122025591797069
    "E13" = 0.1028569265416941
    "G13" = 3.498806644413094
    "H13" = 2.509107105233511
    "K13" = 0.5560565711436993
    "L13" = 0.2102773645743667
    "M13" = 0.2428074885166893
    "B14" = 0.9829481051103812
    "C14" = 0.09443692240949986
    "D14" = 0.1106235095143688
    "E14" = 0.1027512752741693
    "G14" = 3.481046546507343
    "H14" = 2.501251060823392
    "K14" = 0.5507910919372989
    "L14" = 0.2096266328601502
    "M14" = 0.2416463918237568
    "B15" = 0.9794103143455857
    "C15" = 0.09417691269890582
    "D15" = 0.1096571584982371
    "E15" = 0.1026869644273241
    "G15" = 3.470182807695835
    "H15" = 2.496449009900573
    "K15" = 0.5475744584963138
    "L15" = 0.2092299599832756
    "M15" = 0.2409381731299192
    "B16" = 0.9593451466005547
    "C16" = 0.09267523327613247
    "D16" = 0.1041325436702323
    "E16" = 0.10232461363886
    "G16" = 3.408153727367193
    "H16" = 2.469083436499659
    "K16" = 0.5292732186507294
    "L16" = 0.2069862033918355
    "M16" = 0.2369254229892377
    "B17" = 0.9472217475295679
    "C17" = 0.09174351815217108
    "D17" = 0.1007550341184071
    "E17" = 0.1021078481946684
    "G17" = 3.370302295414518
    "H17" = 2.452431798238109
    "K17" = 0.5181635526546131
    "L17" = 0.2056359804393537
    "M17" = 0.2345045562117463
    "B18" = 0.9403169796598263
    "C18" = 0.09120373948398708
    "D18" = 0.09881663704169341
    "E18" = 0.1019852019801952
    "G18" = 3.348604845174776
    "H18" = 2.442904169049768
    "K18" = 0.5118166795369063
    "L18" = 0.2048690149311625
    "M18" = 0.2331271365843008
    "B19" = 0.9379908693123014
    "C19" = 0.09102031050984749
    "D19" = 0.09816105752386761
    "E19" = 0.1019440251033732
    "G19" = 3.341271087640791
    "H19" = 2.439686845552814
    "K19" = 0.5096751338668071
    "L19" = 0.2046109905384697
    "M19" = 0.232663341652426
    "B20" = 0.9485052353931565
    "C20" = 0.09184310182681799
    "D20" = 0.1011141345305049
    "E20" = 0.1021307130150007
    "G20" = 3.374324010869742
    "H20" = 2.454199220965904
    "K20" = 0.5193417315976774
    "L20" = 0.2057787155619621
    "M20" = 0.2347607091900485
    "B21" = 0.9846482612263969
    "C21" = 0.09456139871860358
    "D21" = 0.111087134771779
    "E21" = 0.1027822235075639
    "G21" = 3.486260024589853
    "H21" = 2.503556486182049
    "K21" = 0.5523358929606559
    "L21" = 0.2098173686043481
    "M21" = 0.2419868116453543
    "B22" = 1.008878778496353
    "C22" = 0.09630400363661096
    "D22" = 0.1176435923817252
    "E22" = 0.103226090632365
    "G22" = 3.560079747605528
    "H22" = 2.53626170345143
    "K22" = 0.574285266375739
    "L22" = 0.2125427570573351
    "M22" = 0.2468431465960919
    "B23" = 0.995890733027494
    "C23" = 0.09537702083150634
    "D23" = 0.114140732445037
    "E23" = 0.1029875383315897
    "G23" = 3.520619768587665
    "H23" = 2.518765252043465
    "K23" = 0.5625350818338006
    "L23" = 0.211080306857184
    "M23" = 0.2442389944254089
    "B24" = 0.9479247676514149
    "C24" = 0.09179809289810947
    "D24" = 0.1009517747197464
    "E24" = 0.1021203696775235
    "G24" = 3.372505594527041
    "H24" = 2.453400027066948
    "K24" = 0.5188089516341563
    "L24" = 0.2057141560589031
    "M24" = 0.2346448576728406
    "B25" = 0.8987959354784607
    "C25" = 0.08778890568191855
    "D25" = 0.08688769153215503
    "E25" = 0.101262667379082
    "G25" = 3.21555322816414
    "H25" = 2.384801053089177
    "K25" = 0.4732902875205127
    "L25" = 0.2002947040920944
    "M25" = 0.2248693759851577
}

foreach ($key in $data.Keys) {
    $ws.Range($key).Value = $data[$key]
}
